# "Update countries & provincias Spain"
#
# The source data for several countries in the "Pais" sheet was refreshed:
#   - Singapur   (row 29)
#   - Rumania    (row 37)
#   - Armenia    (row 69)
#   - Cuba / Bulgaria (rows 79/80 - these two countries also swapped places,
#     i.e. the row that used to show "Cuba" now shows "Bulgaria" and vice
#     versa, each with its own refreshed figures)
#   - Georgia    (row 110)
#   - Montenegro (row 128)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 29: Singapur ---------------------------------------------------
$ws.Range("B29").Value = 18205
$ws.Range("C29").Value = 657
$ws.Range("E29").Value = 16841

# --- Row 37: Rumania ------------------------------------------------------
$ws.Range("E37").Value = 7405
$ws.Range("G37").Value = 9
$ws.Range("H37").Value = 780

# --- Row 69: Armenia --------------------------------------------------
$ws.Range("B69").Value = 2386
$ws.Range("C69").Value = 113
$ws.Range("D69").Value = 1035
$ws.Range("E69").Value = 1316
$ws.Range("G69").Value = 2
$ws.Range("H69").Value = 35

# --- Rows 79/80: Cuba and Bulgaria swap places, with updated figures ----
# Row 79 used to be Cuba; it now shows Bulgaria's (refreshed) data.
$ws.Range("A79").Value = "Bulgaria"
$ws.Range("B79").Value = 1611
$ws.Range("C79").Value = 17
$ws.Range("D79").Value = 308
$ws.Range("E79").Value = 1231
$ws.Range("F79").Value = 40
$ws.Range("G79").Value = 0
$ws.Range("H79").Value = 72

# Row 80 used to be Bulgaria; it now shows Cuba's (refreshed) data.
$ws.Range("A80").Value = "Cuba"
$ws.Range("B80").Value = 1611
$ws.Range("C80").Value = 0
$ws.Range("D80").Value = 765
$ws.Range("E80").Value = 780
$ws.Range("F80").Value = 10
$ws.Range("G80").Value = 0
$ws.Range("H80").Value = 66

# --- Row 110: Georgia ---------------------------------------------------
$ws.Range("B110").Value = 589
$ws.Range("C110").Value = 7
$ws.Range("D110").Value = 215
$ws.Range("E110").Value = 366

# --- Row 128: Montenegro -------------------------------------------------
$ws.Range("D128").Value = 249
$ws.Range("E128").Value = 65
